$d = $word.ActiveDocument

$replacements = @(
    @{old = "174÷9="; new = "556÷7="},
    @{old = "183÷5="; new = "182÷7="},
    @{old = "522÷9="; new = "504÷9="},
    @{old = "797÷7="; new = "775÷9="},
    @{old = "212÷8="; new = "511÷3="},
    @{old = "305÷4="; new = "749÷5="},
    @{old = "581÷2="; new = "988÷5="},
    @{old = "683÷2="; new = "128÷9="},
    @{old = "153÷8="; new = "390÷7="},
    @{old = "509÷9="; new = "131÷9="},
    @{old = "335÷2="; new = "459÷7="},
    @{old = "697÷4="; new = "974÷6="},
    @{old = "118÷2="; new = "621÷2="},
    @{old = "658÷5="; new = "275÷8="},
    @{old = "154÷5="; new = "167÷8="},
    @{old = "478÷6="; new = "798÷7="},
    @{old = "123÷4="; new = "922÷8="},
    @{old = "701÷7="; new = "886÷2="},
    @{old = "415÷7="; new = "894÷9="},
    @{old = "664÷9="; new = "436÷2="},
    @{old = "848÷8="; new = "374÷3="},
    @{old = "458÷8="; new = "931÷3="},
    @{old = "212÷2="; new = "541÷3="},
    @{old = "631÷2="; new = "158÷8="},
    @{old = "456÷2="; new = "203÷8="}
)

foreach ($r in $replacements) {
    $d.Content.Find.Execute($r.old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $r.new, 2)
}
